$wb = $excel.ActiveWorkbook

# Hyperlink targets (same "source .md on GitHub" targets used by the existing A2/A3 links)
$url7a8 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef447eefbac355a648b6f335cc5f07434b88f3c0/e2e/7a825a0a-91bc-4596-b72c-2c63a4a27d1c.md"
$urlae0 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef447eefbac355a648b6f335cc5f07434b88f3c0/e2e/ae04b1a7-1635-461b-a462-4610522695de.md"

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status went from "Ready for handoff" to "Handed back: in sync with en-US" ---
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("E2").Value = $newStatus
$wsOv.Range("F2").Value = $newStatus
$wsOv.Range("E3").Value = $newStatus
$wsOv.Range("F3").Value = $newStatus
$wsOv.Columns.Item(5).ColumnWidth = 29.14
$wsOv.Columns.Item(6).ColumnWidth = 29.14

# --- zh-cn sheet: handback for both files recorded ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $url7a8, "", "", "7a825a0a-91bc-4596-b72c-2c63a4a27d1c.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlae0, "", "", "ae04b1a7-1635-461b-a462-4610522695de.md")
$wsZh.Range("J2").Value = "7a825a0a-91bc-4596-b72c-2c63a4a27d1c.2f97a3d5f79622a1251b639e90279f9a1342dc06.zh-cn.xlf"
$wsZh.Range("J3").Value = "ae04b1a7-1635-461b-a462-4610522695de.1476d98eaf2805d158df598827d808fb07eab00e.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-18 07:01:30"
$wsZh.Range("K3").Value = "2016-10-18 07:01:30"
$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# --- de-de sheet: handback for both files recorded ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $url7a8, "", "", "7a825a0a-91bc-4596-b72c-2c63a4a27d1c.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlae0, "", "", "ae04b1a7-1635-461b-a462-4610522695de.md")
$wsDe.Range("J2").Value = "7a825a0a-91bc-4596-b72c-2c63a4a27d1c.2f97a3d5f79622a1251b639e90279f9a1342dc06.de-de.xlf"
$wsDe.Range("J3").Value = "ae04b1a7-1635-461b-a462-4610522695de.1476d98eaf2805d158df598827d808fb07eab00e.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-18 07:02:02"
$wsDe.Range("K3").Value = "2016-10-18 07:02:02"
$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
